$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new week's record at row 131, shifting the
# existing rows 131:241 down to 132:242 (dimension grows to A1:T242).
$ws.Rows.Item(131).Insert()

$ws.Range("A131").Value = 5
$ws.Range("B131").Value = 'Macroferia Regional de Talca'
$ws.Range("C131").Value = 'Maule'
$ws.Range("D131").Value = 44669
$ws.Range("E131").Value = 7
$ws.Range("F131").Value = 'Fruta'
$ws.Range("G131").Value = 100108
$ws.Range("H131").Value = 'Tropicales y subtropicales'
$ws.Range("I131").Value = 100108005
$ws.Range("J131").Value = 'Piña'
$ws.Range("K131").Value = 'Caramelo'
$ws.Range("L131").Value = 'Segunda'
$ws.Range("M131").Value = 540
$ws.Range("N131").Value = 15000
$ws.Range("O131").Value = 15000
$ws.Range("P131").Value = 15000
$ws.Range("Q131").Value = '$/caja 14 unidades'
$ws.Range("R131").Value = 'Ecuador'
$ws.Range("S131").Value = 1071
$ws.Range("T131").Value = 14
